# Update "想去人数" (column F) values on both the "展览" and "全部类型"
# worksheets, which hold duplicate copies of the same convention list.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    6  = 571
    8  = 2032
    11 = 4400
    16 = 116
    17 = 27
    20 = 3223
    22 = 478
    25 = 76
    29 = 57
    32 = 578
    33 = 1844
    34 = 285
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
